$wb = $excel.ActiveWorkbook

# --- Clear (undo) the January attendance marks for AU:BC on rows 6-8 of "Registros" ---
$ws = $wb.Worksheets.Item("Registros")

# Row 6 (joão do pão): clear AU6:AX6, keep AY6, clear AZ6:BC6
$ws.Range("AU6:AX6").Value = ""
$ws.Range("AZ6:BC6").Value = ""

# Row 7 (fernando lando): clear AU7:AX7, keep AY7, clear AZ7:BC7
$ws.Range("AU7:AX7").Value = ""
$ws.Range("AZ7:BC7").Value = ""

# Row 8 (ana cintra): clear AU8:AX8, change AY8 from "c" to "j", clear AZ8:BC8
$ws.Range("AU8:AX8").Value = ""
$ws.Range("AY8").Value = "j"
$ws.Range("AZ8:BC8").Value = ""

# --- Add the new "Justificativas" sheet with the justification records ---
# Insert it right after "Registros" so sheet order matches the target workbook
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "Justificativas"

$newSheet.Range("A1").Value = "Nome"
$newSheet.Range("B1").Value = "Data"
$newSheet.Range("C1").Value = "Motivo"

# Match the bold/centered/bordered header style used by the other sheets
$ws.Range("A1:C1").Copy()
$newSheet.Range("A1:C1").PasteSpecial(-4122)

$newSheet.Range("A2").Value = "ana cintra"
$newSheet.Range("B2").Value = "23/01/2026"
$newSheet.Range("C2").Value = "médico"

$newSheet.Range("A3").Value = "ana cintra"
$newSheet.Range("B3").Value = "30/01/2026"
$newSheet.Range("C3").Value = "frio"
